$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.075.14'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.896.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.50%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '465.76'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +8.66%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.19'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.29%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000345'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.02'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.13%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.521.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.07%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.37'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.870.48'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.38%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.95'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.37%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.82%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.233.96'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.32%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.65'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.36%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.34'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.00%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.66'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.56'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.52'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.77%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.70'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.47%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.60'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.77%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '739.86'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.80%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.60'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.36%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.16%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.48%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.13'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.158'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.17%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.97'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.66%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.08%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0785'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +13.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.35'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -8.60%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.20'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +12.24%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0476'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.21%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.13%  '

$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.140'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.334'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.52%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.56%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.41'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.15%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.47'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.30%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.16'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.40%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.45'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.39%  '
